$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Stage the 4 distinct cell formats used in the edited range into scratch cells ----
# (captured from their original locations before any destructive edits happen)
$ws.Range("A3").Copy() | Out-Null
$ws.Range("Z1").PasteSpecial(-4122) | Out-Null
$ws.Range("A11").Copy() | Out-Null
$ws.Range("Z2").PasteSpecial(-4122) | Out-Null
$ws.Range("B10").Copy() | Out-Null
$ws.Range("Z3").PasteSpecial(-4122) | Out-Null
$ws.Range("C6").Copy() | Out-Null
$ws.Range("Z4").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# ---- Fully remove cells that no longer exist in the new layout ----
$ws.Range("A11").Clear() | Out-Null
$ws.Range("C12").Clear() | Out-Null
$ws.Range("C13").Clear() | Out-Null
$ws.Range("A16").Clear() | Out-Null
$ws.Range("A21").Clear() | Out-Null
$ws.Range("C25").Clear() | Out-Null
$ws.Range("C26").Clear() | Out-Null

# ---- Apply the correct (pre-existing) style to every target cell ----
$ws.Range("Z1").Copy() | Out-Null
$ws.Range("A3").PasteSpecial(-4122) | Out-Null
$ws.Range("Z1").Copy() | Out-Null
$ws.Range("B4").PasteSpecial(-4122) | Out-Null
$ws.Range("Z1").Copy() | Out-Null
$ws.Range("C4").PasteSpecial(-4122) | Out-Null
$ws.Range("Z1").Copy() | Out-Null
$ws.Range("B5").PasteSpecial(-4122) | Out-Null
$ws.Range("Z1").Copy() | Out-Null
$ws.Range("C5").PasteSpecial(-4122) | Out-Null
$ws.Range("Z1").Copy() | Out-Null
$ws.Range("B6").PasteSpecial(-4122) | Out-Null
$ws.Range("Z4").Copy() | Out-Null
$ws.Range("C6").PasteSpecial(-4122) | Out-Null
$ws.Range("Z1").Copy() | Out-Null
$ws.Range("B7").PasteSpecial(-4122) | Out-Null
$ws.Range("Z1").Copy() | Out-Null
$ws.Range("C7").PasteSpecial(-4122) | Out-Null
$ws.Range("Z1").Copy() | Out-Null
$ws.Range("B8").PasteSpecial(-4122) | Out-Null
$ws.Range("Z1").Copy() | Out-Null
$ws.Range("C8").PasteSpecial(-4122) | Out-Null
$ws.Range("Z1").Copy() | Out-Null
$ws.Range("B9").PasteSpecial(-4122) | Out-Null
$ws.Range("Z4").Copy() | Out-Null
$ws.Range("C9").PasteSpecial(-4122) | Out-Null
$ws.Range("Z3").Copy() | Out-Null
$ws.Range("B10").PasteSpecial(-4122) | Out-Null
$ws.Range("Z1").Copy() | Out-Null
$ws.Range("C10").PasteSpecial(-4122) | Out-Null
$ws.Range("Z3").Copy() | Out-Null
$ws.Range("B11").PasteSpecial(-4122) | Out-Null
$ws.Range("Z1").Copy() | Out-Null
$ws.Range("C11").PasteSpecial(-4122) | Out-Null
$ws.Range("Z3").Copy() | Out-Null
$ws.Range("B12").PasteSpecial(-4122) | Out-Null
$ws.Range("Z2").Copy() | Out-Null
$ws.Range("A13").PasteSpecial(-4122) | Out-Null
$ws.Range("Z3").Copy() | Out-Null
$ws.Range("B13").PasteSpecial(-4122) | Out-Null
$ws.Range("Z3").Copy() | Out-Null
$ws.Range("B14").PasteSpecial(-4122) | Out-Null
$ws.Range("Z1").Copy() | Out-Null
$ws.Range("C14").PasteSpecial(-4122) | Out-Null
$ws.Range("Z3").Copy() | Out-Null
$ws.Range("B15").PasteSpecial(-4122) | Out-Null
$ws.Range("Z1").Copy() | Out-Null
$ws.Range("C15").PasteSpecial(-4122) | Out-Null
$ws.Range("Z3").Copy() | Out-Null
$ws.Range("B16").PasteSpecial(-4122) | Out-Null
$ws.Range("Z1").Copy() | Out-Null
$ws.Range("C16").PasteSpecial(-4122) | Out-Null
$ws.Range("Z3").Copy() | Out-Null
$ws.Range("B17").PasteSpecial(-4122) | Out-Null
$ws.Range("Z1").Copy() | Out-Null
$ws.Range("C17").PasteSpecial(-4122) | Out-Null
$ws.Range("Z3").Copy() | Out-Null
$ws.Range("B18").PasteSpecial(-4122) | Out-Null
$ws.Range("Z4").Copy() | Out-Null
$ws.Range("C18").PasteSpecial(-4122) | Out-Null
$ws.Range("Z3").Copy() | Out-Null
$ws.Range("B19").PasteSpecial(-4122) | Out-Null
$ws.Range("Z1").Copy() | Out-Null
$ws.Range("C19").PasteSpecial(-4122) | Out-Null
$ws.Range("Z3").Copy() | Out-Null
$ws.Range("B20").PasteSpecial(-4122) | Out-Null
$ws.Range("Z1").Copy() | Out-Null
$ws.Range("C20").PasteSpecial(-4122) | Out-Null
$ws.Range("Z3").Copy() | Out-Null
$ws.Range("B21").PasteSpecial(-4122) | Out-Null
$ws.Range("Z1").Copy() | Out-Null
$ws.Range("C21").PasteSpecial(-4122) | Out-Null
$ws.Range("Z3").Copy() | Out-Null
$ws.Range("B22").PasteSpecial(-4122) | Out-Null
$ws.Range("Z1").Copy() | Out-Null
$ws.Range("C22").PasteSpecial(-4122) | Out-Null
$ws.Range("Z3").Copy() | Out-Null
$ws.Range("B23").PasteSpecial(-4122) | Out-Null
$ws.Range("Z1").Copy() | Out-Null
$ws.Range("C23").PasteSpecial(-4122) | Out-Null
$ws.Range("Z3").Copy() | Out-Null
$ws.Range("B24").PasteSpecial(-4122) | Out-Null
$ws.Range("Z1").Copy() | Out-Null
$ws.Range("C24").PasteSpecial(-4122) | Out-Null
$ws.Range("Z3").Copy() | Out-Null
$ws.Range("B25").PasteSpecial(-4122) | Out-Null
$ws.Range("Z2").Copy() | Out-Null
$ws.Range("A26").PasteSpecial(-4122) | Out-Null
$ws.Range("Z3").Copy() | Out-Null
$ws.Range("B26").PasteSpecial(-4122) | Out-Null
$ws.Range("Z3").Copy() | Out-Null
$ws.Range("B27").PasteSpecial(-4122) | Out-Null
$ws.Range("Z1").Copy() | Out-Null
$ws.Range("C27").PasteSpecial(-4122) | Out-Null
$ws.Range("Z3").Copy() | Out-Null
$ws.Range("B28").PasteSpecial(-4122) | Out-Null
$ws.Range("Z1").Copy() | Out-Null
$ws.Range("C28").PasteSpecial(-4122) | Out-Null
$ws.Range("Z3").Copy() | Out-Null
$ws.Range("B29").PasteSpecial(-4122) | Out-Null
$ws.Range("Z4").Copy() | Out-Null
$ws.Range("C29").PasteSpecial(-4122) | Out-Null
$ws.Range("Z3").Copy() | Out-Null
$ws.Range("B30").PasteSpecial(-4122) | Out-Null
$ws.Range("Z1").Copy() | Out-Null
$ws.Range("C30").PasteSpecial(-4122) | Out-Null
$ws.Range("Z3").Copy() | Out-Null
$ws.Range("B31").PasteSpecial(-4122) | Out-Null
$ws.Range("Z4").Copy() | Out-Null
$ws.Range("C31").PasteSpecial(-4122) | Out-Null
$ws.Range("Z3").Copy() | Out-Null
$ws.Range("B32").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# ---- Make sure previously-numbered rows that are now blank truly have no value ----
$ws.Range("B12").ClearContents() | Out-Null
$ws.Range("B13").ClearContents() | Out-Null
$ws.Range("B25").ClearContents() | Out-Null
$ws.Range("B26").ClearContents() | Out-Null

# ---- Write the new cell values (numbers stay numeric, text stays text) ----
$values = @{
    "A3" = 62.1
    "B4" = 2
    "C4" = 'doing'
    "B5" = 3
    "C5" = 'coming'
    "B6" = 4
    "C6" = 'going(spending/having)'
    "B7" = 5
    "C7" = 'buying'
    "B8" = 6
    "C8" = 'seeing'
    "B9" = 7
    "C9" = 'playing(watching)'
    "B10" = 8
    "C10" = 'solving'
    "B11" = 9
    "C11" = 'having'
    "A13" = '62.2'
    "B14" = 2
    "C14" = 'of causing'
    "B15" = 3
    "C15" = 'from walking'
    "B16" = 4
    "C16" = 'for interrupting'
    "B17" = 5
    "C17" = 'of using'
    "B18" = 6
    "C18" = 'of/about doing'
    "B19" = 7
    "C19" = 'from escaping'
    "B20" = 8
    "C20" = 'on telling'
    "B21" = 9
    "C21" = 'to eating'
    "B22" = 10
    "C22" = 'for being'
    "B23" = 11
    "C23" = 'for inviting'
    "B24" = 12
    "C24" = 'of,wearing'
    "A26" = '62.3'
    "B27" = 2
    "C27" = 'on taking to the station'
    "B28" = 3
    "C28" = 'on getting married'
    "B29" = 4
    "C29" = 'me(Sue) for coming to see her'
    "B30" = 5
    "C30" = 'to me for not phoning earlier'
    "B31" = 6
    "C31" = 'me of (being) selfish'
}
foreach ($cellRef in $values.Keys) {
    $ws.Range($cellRef).Value = $values[$cellRef]
}

# ---- Remove the scratch cells used for style staging ----
$ws.Range("Z1:Z4").Clear() | Out-Null

# ---- Rich text formatting: strike through "/about" in the C18 answer ----
$ws.Range("C18").Characters(3, 6).Font.Strikethrough = $true

# ---- Restore the view/selection state recorded in the saved workbook ----
$excel.Goto($ws.Range("A25"), $true)
$ws.Range("C29").Select() | Out-Null
